$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column H (Datum) to text format so date-like strings are not
# auto-converted to Excel date serials.
$ws.Range("H2:H21").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = 20
$ws.Range("B2").Value = "ORD020"
$ws.Range("C2").Value = "Лукас Крумпах"
$ws.Range("D2").Value = "lukas.krumpach@gmail.com"
$ws.Range("E2").Value = "Moscow, Потаповский переулок, д. 8/12 стр. 2"
$ws.Range("F2").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = "2025-03-26"
$ws.Range("I2").Value = "F020"
$ws.Range("J2").Value = "Nová"

# Row 3
$ws.Range("A3").Value = 19
$ws.Range("B3").Value = "ORD019"
$ws.Range("C3").Value = "Лукас Крумпах"
$ws.Range("D3").Value = "lukas.krumpach@gmail.com"
$ws.Range("E3").Value = "Moscow, Потаповский переулок, д. 8/12 стр. 2"
$ws.Range("F3").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = "2025-03-25"
$ws.Range("I3").Value = "F019"
$ws.Range("J3").Value = "Zpracovává se"

# Row 4
$ws.Range("A4").Value = 18
$ws.Range("B4").Value = "ORD018"
$ws.Range("C4").Value = "Hana VAVROVA"
$ws.Range("D4").Value = "lukas.krumpach@gmail.com"
$ws.Range("E4").Value = "Jaktáře 14756666666666666666666666"
$ws.Range("F4").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = "2025-03-25"
$ws.Range("I4").Value = "F018"
$ws.Range("J4").Value = "Nová"

# Row 5
$ws.Range("A5").Value = 17
$ws.Range("B5").Value = "ORD017"
$ws.Range("C5").Value = "Lukas Krumpach"
$ws.Range("D5").Value = "lukas.krumpach@gmail.com"
$ws.Range("E5").Value = "Pod Lipami`nC.P. 24 44444444444444444444"
$ws.Range("F5").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = "2025-03-25"
$ws.Range("I5").Value = "F017"
$ws.Range("J5").Value = "Nová"

# Row 6
$ws.Range("A6").Value = 16
$ws.Range("B6").Value = "ORD016"
$ws.Range("C6").Value = "Hana VAVROVA"
$ws.Range("D6").Value = "lukas.krumpach@gmail.com"
$ws.Range("E6").Value = "Jaktáře 1475 3333333333333333333"
$ws.Range("F6").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = "2025-03-25"
$ws.Range("I6").Value = "F016"
$ws.Range("J6").Value = "Nová"

# Row 7
$ws.Range("A7").Value = 15
$ws.Range("B7").Value = "ORD015"
$ws.Range("C7").Value = "Lukas Krumpach"
$ws.Range("D7").Value = "lukas.krumpach@gmail.com"
$ws.Range("E7").Value = "Premonstrátů 910C 22222222222222222222222222222"
$ws.Range("F7").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = "2025-03-25"
$ws.Range("I7").Value = "F015"
$ws.Range("J7").Value = "Odesláno"

# Row 8
$ws.Range("A8").Value = 14
$ws.Range("B8").Value = "ORD014"
$ws.Range("C8").Value = "Hana VAVROVA"
$ws.Range("D8").Value = "lukas.krumpach@gmail.com"
$ws.Range("E8").Value = "Jaktáře 1475"
$ws.Range("F8").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = "2025-03-25"
$ws.Range("I8").Value = "F014"
$ws.Range("J8").Value = "Nová"

# Row 9
$ws.Range("A9").Value = 13
$ws.Range("B9").Value = "ORD013"
$ws.Range("C9").Value = "Лукас Крумпах"
$ws.Range("D9").Value = "lukas.krumpach@gmail.com"
$ws.Range("E9").Value = "Moscow, Потаповский переулок, д. 8/12 стр. 2"
$ws.Range("F9").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = "2025-03-25"
$ws.Range("I9").Value = "F013"
$ws.Range("J9").Value = "Zpracovává se"

# Row 10
$ws.Range("A10").Value = 12
$ws.Range("B10").Value = "ORD012"
$ws.Range("C10").Value = "Mariia Isova"
$ws.Range("D10").Value = "l.m.e.companycz@gmail.com"
$ws.Range("E10").Value = "Názovská 14"
$ws.Range("F10").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = "2025-03-24"
$ws.Range("I10").Value = "F012"
$ws.Range("J10").Value = "Nová"

# Row 11
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "ORD011"
$ws.Range("C11").Value = "Лукас Крумпах"
$ws.Range("D11").Value = "lukas.krumpach@gmail.com"
$ws.Range("E11").Value = "Moscow, Потаповский переулок, д. 8/12 стр. 2"
$ws.Range("F11").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = "2025-03-24"
$ws.Range("I11").Value = "F011"
$ws.Range("J11").Value = "Nová"

# Row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "ORD010"
$ws.Range("C12").Value = "Lukas Krumpach"
$ws.Range("D12").Value = "lukas.krumpach@gmail.com"
$ws.Range("E12").Value = "Premonstrátů 910C"
$ws.Range("F12").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = "2025-03-24"
$ws.Range("I12").Value = "F010"
$ws.Range("J12").Value = "Nová"

# Row 13
$ws.Range("A13").Value = 9
$ws.Range("B13").Value = "ORD009"
$ws.Range("C13").Value = "Lukas Krumpach"
$ws.Range("D13").Value = "lukas.krumpach@gmail.com"
$ws.Range("E13").Value = "Premonstrátů 910C"
$ws.Range("F13").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = "2025-03-24"
$ws.Range("I13").Value = "F009"
$ws.Range("J13").Value = "Nová"

# Row 14
$ws.Range("A14").Value = 8
$ws.Range("B14").Value = "ORD008"
$ws.Range("C14").Value = "Lukas Krumpach"
$ws.Range("D14").Value = "lukas.krumpach@gmail.com"
$ws.Range("E14").Value = "Premonstrátů 910C"
$ws.Range("F14").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = "2025-03-23"
$ws.Range("I14").Value = "F008"
$ws.Range("J14").Value = "Nová"

# Row 15
$ws.Range("A15").Value = 7
$ws.Range("B15").Value = "ORD007"
$ws.Range("C15").Value = "Lukas Krumpach"
$ws.Range("D15").Value = "lukas.krumpach@gmail.com"
$ws.Range("E15").Value = "Premonstrátů 910C"
$ws.Range("F15").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = "2025-03-23"
$ws.Range("I15").Value = "F007"
$ws.Range("J15").Value = "Nová"

# Row 16
$ws.Range("A16").Value = 6
$ws.Range("B16").Value = "ORD006"
$ws.Range("C16").Value = "Lukas Krumpach"
$ws.Range("D16").Value = "lukas.krumpach@gmail.com"
$ws.Range("E16").Value = "Premonstrátů 910C"
$ws.Range("F16").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = "2025-03-23"
$ws.Range("I16").Value = "F006"
$ws.Range("J16").Value = "Nová"

# Row 17
$ws.Range("A17").Value = 5
$ws.Range("B17").Value = "ORD005"
$ws.Range("C17").Value = "Лукас Крумпах"
$ws.Range("D17").Value = "lukas.krumpach@gmail.com"
$ws.Range("E17").Value = "Moscow, Потаповский переулок, д. 8/12 стр. 2"
$ws.Range("F17").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = "2025-03-23"
$ws.Range("I17").Value = "F005"
$ws.Range("J17").Value = "Nová"

# Row 18
$ws.Range("A18").Value = 4
$ws.Range("B18").Value = "ORD004"
$ws.Range("C18").Value = "Lukas Krumpach"
$ws.Range("D18").Value = "lukas.krumpach@gmail.com"
$ws.Range("E18").Value = "Premonstrátů 910C"
$ws.Range("F18").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = "2025-03-23"
$ws.Range("I18").Value = "F004"
$ws.Range("J18").Value = "Nová"

# Row 19
$ws.Range("A19").Value = 3
$ws.Range("B19").Value = "ORD003"
$ws.Range("C19").Value = "Lukas Krumpach"
$ws.Range("D19").Value = "lukas.krumpach@gmail.com"
$ws.Range("E19").Value = "Premonstrátů 910C"
$ws.Range("F19").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = "2025-03-23"
$ws.Range("I19").Value = "F003"
$ws.Range("J19").Value = "Nová"

# Row 20
$ws.Range("A20").Value = 2
$ws.Range("B20").Value = "ORD002"
$ws.Range("C20").Value = "Лукас Крумпах"
$ws.Range("D20").Value = "lukas.krumpach@gmail.com"
$ws.Range("E20").Value = "Moscow, Потаповский переулок, д. 8/12 стр. 2"
$ws.Range("F20").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = "2025-03-23"
$ws.Range("I20").Value = "F002"
$ws.Range("J20").Value = "Nová"

# Row 21
$ws.Range("A21").Value = 1
$ws.Range("B21").Value = "ORD001"
$ws.Range("C21").Value = "Lukas Krumpach"
$ws.Range("D21").Value = "lukas.krumpach@gmail.com"
$ws.Range("E21").Value = "Premonstrátů 910C"
$ws.Range("F21").Value = "🔥 HYALCHONDRO® HC PLUS"
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = "2025-03-23"
$ws.Range("I21").Value = "F001"
$ws.Range("J21").Value = "Nová"
